$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-27 Thursday" "2024-06-28 Friday"

Replace-Text "175÷9=" "288÷9="
Replace-Text "480÷3=" "898÷8="
Replace-Text "272÷6=" "424÷4="
Replace-Text "531÷5=" "130÷7="
Replace-Text "928÷2=" "230÷7="

Replace-Text "782÷8=" "355÷2="
Replace-Text "663÷4=" "361÷8="
Replace-Text "417÷5=" "269÷4="
Replace-Text "588÷7=" "552÷4="
Replace-Text "814÷8=" "584÷8="

Replace-Text "797÷9=" "310÷8="
Replace-Text "545÷8=" "530÷3="
Replace-Text "147÷4=" "674÷6="
Replace-Text "227÷8=" "245÷3="
Replace-Text "891÷7=" "699÷4="

Replace-Text "764÷2=" "640÷3="
Replace-Text "650÷7=" "926÷2="
Replace-Text "903÷3=" "363÷9="
Replace-Text "307÷6=" "286÷4="
Replace-Text "680÷7=" "285÷6="

Replace-Text "641÷4=" "670÷9="
Replace-Text "706÷8=" "182÷2="
Replace-Text "212÷3=" "454÷2="
Replace-Text "425÷9=" "131÷7="
Replace-Text "249÷5=" "912÷5="
